$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rows 308-313 (columns E, F, G changed)
$ws.Range("E308").Value = 1574
$ws.Range("F308").Value = 5053
$ws.Range("G308").Value = 80

$ws.Range("E309").Value = 726
$ws.Range("F309").Value = 3088
$ws.Range("G309").Value = 56

$ws.Range("E310").Value = 767
$ws.Range("F310").Value = 3214
$ws.Range("G310").Value = 68

$ws.Range("E311").Value = 1403
$ws.Range("F311").Value = 10532
$ws.Range("G311").Value = 74

$ws.Range("E312").Value = 457
$ws.Range("F312").Value = 3409
$ws.Range("G312").Value = 20

$ws.Range("E313").Value = 299
$ws.Range("F313").Value = 3027
$ws.Range("G313").Value = 25

# Add new row 314 (A314 must stay plain text "11.01.2021", not an
# auto-converted date serial; force text format, write, then drop the
# format override so the cell ends up with no explicit style, matching
# the other date cells in column A)
$ws.Range("A314").NumberFormat = "@"
$ws.Range("A314").Value = "11.01.2021"
$ws.Range("A314").ClearFormats()
$ws.Range("B314").Value = 115633
$ws.Range("C314").Value = 720
$ws.Range("D314").Value = 3280815
$ws.Range("E314").Value = 567
$ws.Range("F314").Value = 2557
$ws.Range("G314").Value = 28
